$wb = $excel.ActiveWorkbook

# --- Sheet: "High Priority break-up" (sheet5) -> rename and update values ---
$wsOld = $wb.Worksheets.Item("High Priority break-up")
$wsOld.Name = "Interannual update - High Pri"

# Update Interannual update - High Pri values
$wsOld.Range("B2").Value = 73
$wsOld.Range("C2").Value = 70.90000000000001
$wsOld.Range("D2").Value = 73
$wsOld.Range("E2").Value = 77.7

$wsOld.Range("B3").Value = 30
$wsOld.Range("C3").Value = 29.1
$wsOld.Range("D3").Value = 21
$wsOld.Range("E3").Value = 22.3

# --- New sheet: "Major update - High Priority " (sheet6) - copy of the original content ---
# Add it after the last existing sheet so it lands at the end of the tab order.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add($null, $lastSheet)
$wsNew.Name = "Major update - High Priority "

$wsNew.Range("A1").Value = "Break-up"
$wsNew.Range("B1").Value = "High Species (no.)"
$wsNew.Range("C1").Value = "High Species (perc.)"
$wsNew.Range("D1").Value = "New High Species (no.)"
$wsNew.Range("E1").Value = "New High Species (perc.)"
$wsNew.Range("A1:E1").Font.Bold = $true
$wsNew.Range("A1:E1").HorizontalAlignment = -4108

$wsNew.Range("A2").Value = "Trend New"
$wsNew.Range("B2").Value = 1
$wsNew.Range("C2").Value = 9.1
$wsNew.Range("D2").Value = 1
$wsNew.Range("E2").Value = 9.1

$wsNew.Range("A3").Value = "IUCN"
$wsNew.Range("B3").Value = 10
$wsNew.Range("C3").Value = 90.90000000000001
$wsNew.Range("D3").Value = 10
$wsNew.Range("E3").Value = 90.90000000000001


# --- Sheet1: "Trends Status" ---
$ws1 = $wb.Worksheets.Item("Trends Status")
$ws1.Range("B2").Value = 0
$ws1.Range("D2").ClearContents()
$ws1.Range("D3").ClearContents()
$ws1.Range("D4").ClearContents()
$ws1.Range("D5").ClearContents()
$ws1.Range("D6").ClearContents()
$ws1.Range("B7").Value = 2
$ws1.Range("C7").Value = 14
$ws1.Range("B8").Value = 384
$ws1.Range("C8").Value = 371

# --- Sheet3: "Priority Status" ---
$ws3 = $wb.Worksheets.Item("Priority Status")
$ws3.Range("B2").Value = 103
$ws3.Range("B3").Value = 286
$ws3.Range("B4").Value = 554

# --- Sheet4: "Species qualification" ---
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("A2").Value = "SoIB Assessment"
$ws4.Range("B2").Value = 386
$ws4.Range("C3").Value = 0
$ws4.Range("B4").Value = 15
